# Fix experiment name in example data, and make the corrected sheet the
# active/selected tab (it was previously misnamed and sheet 1 was active).

$wb = $excel.ActiveWorkbook

# Correct the misspelled / truncated experiment name on the 4th sheet.
$ws = $wb.Worksheets.Item("MSSA Biofilm Planktonic Inhibit")
$ws.Name = "MSSA Planktonic Inhibition"

# Make the renamed sheet the active tab (was previously sheet 1).
$ws.Activate()
